{"js": "// Add a period to the end of the \"Used iframe...\" bullet, then add a new\n// bullet list item after it for the TMDB API work.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Find the paragraph that ends the iframe bullet point.\nconst target = paragraphs.items.find((p) =>\n  p.text.indexOf(\"Used iframe to embed the preview of a movie\") !== -1\n);\nif (!target) {\n  throw new Error(\"Could not find target paragraph\");\n}\n\n// Append the missing period to the existing sentence.\ntarget.insertText(\".\", \"End\");\n\n// Insert a new list paragraph right after it, inheriting the same\n// (ListParagraph / bullet) formatting, with the new TMDB API bullet text.\ntarget.insertParagraph(\"TMDB API to get videos and movies and titles.\", \"After\");\n\nawait context.sync();\n", "ps1": "# Add the missing period to the \"Used iframe...\" bullet, then add a new\n# bullet list item after it describing the new TMDB API work.\n$d = $word.ActiveDocument\n\n# Locate the sentence that needs the trailing period via Find; this\n# collapses $r to the matched text so InsertAfter lands right after it.\n$r = $d.Content\n$found = $r.Find.Execute(\"Used iframe to embed the preview of a movie\")\n\nif ($found) {\n    $r.InsertAfter(\".\")\n\n    # Insert a new paragraph right after the (now-perioded) sentence; it\n    # inherits the ListParagraph / bullet formatting of the source paragraph.\n    $r.InsertParagraphAfter()\n\n    $newPara = $d.Paragraphs.Last\n    $newPara.Range.Text = \"TMDB API to get videos and movies and titles.\"\n}\n"}
